# edit.ps1
# Applies "Update latest output (run 263)" to optimisation_result.xlsx:
#  - Schedule: recompute E3/F3, append new pump-run block as row 4
#  - Detailed: reclassify rows 14-33 forecast->historical w/ revised
#    historical prices (rows 14-48 updated), and append a new day of
#    forecast rows (50-97) for 2026-02-28

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Schedule"
# ---------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Schedule")

# Row 3: recomputed Cost ($) / Unit Cost ($/ML)
$wsSchedule.Cells.Item(3, 5).Value = 85.91514749999997
$wsSchedule.Cells.Item(3, 6).Value = 3.246982142857142

# Row 4 (new): additional pump-run block on 2026-02-28
$wsSchedule.Cells.Item(4, 1).Value = 46081.29166666666
$wsSchedule.Cells.Item(4, 2).Value = 46081.75
$wsSchedule.Range("A4:B4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsSchedule.Cells.Item(4, 3).Value = 11
$wsSchedule.Cells.Item(4, 4).Value = 41.58
$wsSchedule.Cells.Item(4, 5).Value = 233.087205
$wsSchedule.Cells.Item(4, 6).Value = 5.605752886002886

# ---------------------------------------------------------------
# Sheet "Detailed"
# ---------------------------------------------------------------
$wsDetailed = $wb.Worksheets.Item("Detailed")

# Revised Price (B) / Type (C) for existing rows 14-48 (2026-02-27)
$wsDetailed.Cells.Item(14, 2).Value = 65
$wsDetailed.Cells.Item(15, 2).Value = 76.47257
$wsDetailed.Cells.Item(16, 2).Value = 75.25190000000001; $wsDetailed.Cells.Item(16, 3).Value = "historical"
$wsDetailed.Cells.Item(17, 2).Value = 57.05689; $wsDetailed.Cells.Item(17, 3).Value = "historical"
$wsDetailed.Cells.Item(18, 2).Value = 51.14382; $wsDetailed.Cells.Item(18, 3).Value = "historical"
$wsDetailed.Cells.Item(19, 2).Value = 37.89; $wsDetailed.Cells.Item(19, 3).Value = "historical"
$wsDetailed.Cells.Item(20, 2).Value = 35.88; $wsDetailed.Cells.Item(20, 3).Value = "historical"
$wsDetailed.Cells.Item(21, 2).Value = 12.94905; $wsDetailed.Cells.Item(21, 3).Value = "historical"
$wsDetailed.Cells.Item(22, 2).Value = 19.564; $wsDetailed.Cells.Item(22, 3).Value = "historical"
$wsDetailed.Cells.Item(23, 2).Value = 0.73; $wsDetailed.Cells.Item(23, 3).Value = "historical"
$wsDetailed.Cells.Item(24, 2).Value = 0.00003; $wsDetailed.Cells.Item(24, 3).Value = "historical"
$wsDetailed.Cells.Item(25, 3).Value = "historical"
$wsDetailed.Cells.Item(26, 3).Value = "historical"
$wsDetailed.Cells.Item(27, 2).Value = 3.61597; $wsDetailed.Cells.Item(27, 3).Value = "historical"
$wsDetailed.Cells.Item(28, 2).Value = 8.299189999999999; $wsDetailed.Cells.Item(28, 3).Value = "historical"
$wsDetailed.Cells.Item(29, 3).Value = "historical"
$wsDetailed.Cells.Item(30, 2).Value = 0.00003; $wsDetailed.Cells.Item(30, 3).Value = "historical"
$wsDetailed.Cells.Item(31, 2).Value = 8.15189; $wsDetailed.Cells.Item(31, 3).Value = "historical"
$wsDetailed.Cells.Item(32, 2).Value = -3.99192; $wsDetailed.Cells.Item(32, 3).Value = "historical"
$wsDetailed.Cells.Item(33, 3).Value = "historical"
$wsDetailed.Cells.Item(34, 2).Value = 22.00169
$wsDetailed.Cells.Item(35, 2).Value = 35.88
$wsDetailed.Cells.Item(36, 2).Value = 48.47932
$wsDetailed.Cells.Item(37, 2).Value = 37.89
$wsDetailed.Cells.Item(38, 2).Value = 48.22872
$wsDetailed.Cells.Item(39, 2).Value = 60.25303
$wsDetailed.Cells.Item(40, 2).Value = 70.01016
$wsDetailed.Cells.Item(41, 2).Value = 71.4118
$wsDetailed.Cells.Item(42, 2).Value = 71.24460999999999
$wsDetailed.Cells.Item(43, 2).Value = 73.48285
$wsDetailed.Cells.Item(44, 2).Value = 61.49528
$wsDetailed.Cells.Item(45, 2).Value = 62.96173
$wsDetailed.Cells.Item(46, 2).Value = 65.11269
$wsDetailed.Cells.Item(47, 2).Value = 65
$wsDetailed.Cells.Item(48, 2).Value = 57.0601

# New rows 50-97: forecast data for 2026-02-28
$wsDetailed.Cells.Item(50, 1).Value = 46081; $wsDetailed.Cells.Item(50, 2).Value = 57.06; $wsDetailed.Cells.Item(50, 3).Value = "forecast"; $wsDetailed.Cells.Item(50, 4).Value = 46081; $wsDetailed.Cells.Item(50, 5).Value = "OFF"
$wsDetailed.Cells.Item(51, 1).Value = 46081.02083333334; $wsDetailed.Cells.Item(51, 2).Value = 57.05998; $wsDetailed.Cells.Item(51, 3).Value = "forecast"; $wsDetailed.Cells.Item(51, 4).Value = 46081; $wsDetailed.Cells.Item(51, 5).Value = "OFF"
$wsDetailed.Cells.Item(52, 1).Value = 46081.04166666666; $wsDetailed.Cells.Item(52, 2).Value = 57.04976; $wsDetailed.Cells.Item(52, 3).Value = "forecast"; $wsDetailed.Cells.Item(52, 4).Value = 46081; $wsDetailed.Cells.Item(52, 5).Value = "OFF"
$wsDetailed.Cells.Item(53, 1).Value = 46081.0625; $wsDetailed.Cells.Item(53, 2).Value = 57.06; $wsDetailed.Cells.Item(53, 3).Value = "forecast"; $wsDetailed.Cells.Item(53, 4).Value = 46081; $wsDetailed.Cells.Item(53, 5).Value = "OFF"
$wsDetailed.Cells.Item(54, 1).Value = 46081.08333333334; $wsDetailed.Cells.Item(54, 2).Value = 57.06; $wsDetailed.Cells.Item(54, 3).Value = "forecast"; $wsDetailed.Cells.Item(54, 4).Value = 46081; $wsDetailed.Cells.Item(54, 5).Value = "OFF"
$wsDetailed.Cells.Item(55, 1).Value = 46081.10416666666; $wsDetailed.Cells.Item(55, 2).Value = 56.98; $wsDetailed.Cells.Item(55, 3).Value = "forecast"; $wsDetailed.Cells.Item(55, 4).Value = 46081; $wsDetailed.Cells.Item(55, 5).Value = "OFF"
$wsDetailed.Cells.Item(56, 1).Value = 46081.125; $wsDetailed.Cells.Item(56, 2).Value = 56.78471; $wsDetailed.Cells.Item(56, 3).Value = "forecast"; $wsDetailed.Cells.Item(56, 4).Value = 46081; $wsDetailed.Cells.Item(56, 5).Value = "OFF"
$wsDetailed.Cells.Item(57, 1).Value = 46081.14583333334; $wsDetailed.Cells.Item(57, 2).Value = 38.65901; $wsDetailed.Cells.Item(57, 3).Value = "forecast"; $wsDetailed.Cells.Item(57, 4).Value = 46081; $wsDetailed.Cells.Item(57, 5).Value = "OFF"
$wsDetailed.Cells.Item(58, 1).Value = 46081.16666666666; $wsDetailed.Cells.Item(58, 2).Value = 38.96592; $wsDetailed.Cells.Item(58, 3).Value = "forecast"; $wsDetailed.Cells.Item(58, 4).Value = 46081; $wsDetailed.Cells.Item(58, 5).Value = "OFF"
$wsDetailed.Cells.Item(59, 1).Value = 46081.1875; $wsDetailed.Cells.Item(59, 2).Value = 39.24123; $wsDetailed.Cells.Item(59, 3).Value = "forecast"; $wsDetailed.Cells.Item(59, 4).Value = 46081; $wsDetailed.Cells.Item(59, 5).Value = "OFF"
$wsDetailed.Cells.Item(60, 1).Value = 46081.20833333334; $wsDetailed.Cells.Item(60, 2).Value = 39.87985; $wsDetailed.Cells.Item(60, 3).Value = "forecast"; $wsDetailed.Cells.Item(60, 4).Value = 46081; $wsDetailed.Cells.Item(60, 5).Value = "OFF"
$wsDetailed.Cells.Item(61, 1).Value = 46081.22916666666; $wsDetailed.Cells.Item(61, 2).Value = 55.83347; $wsDetailed.Cells.Item(61, 3).Value = "forecast"; $wsDetailed.Cells.Item(61, 4).Value = 46081; $wsDetailed.Cells.Item(61, 5).Value = "OFF"
$wsDetailed.Cells.Item(62, 1).Value = 46081.25; $wsDetailed.Cells.Item(62, 2).Value = 57.06; $wsDetailed.Cells.Item(62, 3).Value = "forecast"; $wsDetailed.Cells.Item(62, 4).Value = 46081; $wsDetailed.Cells.Item(62, 5).Value = "OFF"
$wsDetailed.Cells.Item(63, 1).Value = 46081.27083333334; $wsDetailed.Cells.Item(63, 2).Value = 56.98; $wsDetailed.Cells.Item(63, 3).Value = "forecast"; $wsDetailed.Cells.Item(63, 4).Value = 46081; $wsDetailed.Cells.Item(63, 5).Value = "OFF"
$wsDetailed.Cells.Item(64, 1).Value = 46081.29166666666; $wsDetailed.Cells.Item(64, 2).Value = 38.1007; $wsDetailed.Cells.Item(64, 3).Value = "forecast"; $wsDetailed.Cells.Item(64, 4).Value = 46081; $wsDetailed.Cells.Item(64, 5).Value = "ON"
$wsDetailed.Cells.Item(65, 1).Value = 46081.3125; $wsDetailed.Cells.Item(65, 2).Value = 35.88; $wsDetailed.Cells.Item(65, 3).Value = "forecast"; $wsDetailed.Cells.Item(65, 4).Value = 46081; $wsDetailed.Cells.Item(65, 5).Value = "ON"
$wsDetailed.Cells.Item(66, 1).Value = 46081.33333333334; $wsDetailed.Cells.Item(66, 2).Value = 0.7; $wsDetailed.Cells.Item(66, 3).Value = "forecast"; $wsDetailed.Cells.Item(66, 4).Value = 46081; $wsDetailed.Cells.Item(66, 5).Value = "ON"
$wsDetailed.Cells.Item(67, 1).Value = 46081.35416666666; $wsDetailed.Cells.Item(67, 2).Value = 4.93013; $wsDetailed.Cells.Item(67, 3).Value = "forecast"; $wsDetailed.Cells.Item(67, 4).Value = 46081; $wsDetailed.Cells.Item(67, 5).Value = "ON"
$wsDetailed.Cells.Item(68, 1).Value = 46081.375; $wsDetailed.Cells.Item(68, 2).Value = 9.07302; $wsDetailed.Cells.Item(68, 3).Value = "forecast"; $wsDetailed.Cells.Item(68, 4).Value = 46081; $wsDetailed.Cells.Item(68, 5).Value = "ON"
$wsDetailed.Cells.Item(69, 1).Value = 46081.39583333334; $wsDetailed.Cells.Item(69, 2).Value = 7.9991; $wsDetailed.Cells.Item(69, 3).Value = "forecast"; $wsDetailed.Cells.Item(69, 4).Value = 46081; $wsDetailed.Cells.Item(69, 5).Value = "ON"
$wsDetailed.Cells.Item(70, 1).Value = 46081.41666666666; $wsDetailed.Cells.Item(70, 2).Value = 11.27986; $wsDetailed.Cells.Item(70, 3).Value = "forecast"; $wsDetailed.Cells.Item(70, 4).Value = 46081; $wsDetailed.Cells.Item(70, 5).Value = "ON"
$wsDetailed.Cells.Item(71, 1).Value = 46081.4375; $wsDetailed.Cells.Item(71, 2).Value = 7.93696; $wsDetailed.Cells.Item(71, 3).Value = "forecast"; $wsDetailed.Cells.Item(71, 4).Value = 46081; $wsDetailed.Cells.Item(71, 5).Value = "ON"
$wsDetailed.Cells.Item(72, 1).Value = 46081.45833333334; $wsDetailed.Cells.Item(72, 2).Value = 0.7; $wsDetailed.Cells.Item(72, 3).Value = "forecast"; $wsDetailed.Cells.Item(72, 4).Value = 46081; $wsDetailed.Cells.Item(72, 5).Value = "ON"
$wsDetailed.Cells.Item(73, 1).Value = 46081.47916666666; $wsDetailed.Cells.Item(73, 2).Value = 0.7; $wsDetailed.Cells.Item(73, 3).Value = "forecast"; $wsDetailed.Cells.Item(73, 4).Value = 46081; $wsDetailed.Cells.Item(73, 5).Value = "ON"
$wsDetailed.Cells.Item(74, 1).Value = 46081.5; $wsDetailed.Cells.Item(74, 2).Value = 0.7; $wsDetailed.Cells.Item(74, 3).Value = "forecast"; $wsDetailed.Cells.Item(74, 4).Value = 46081; $wsDetailed.Cells.Item(74, 5).Value = "ON"
$wsDetailed.Cells.Item(75, 1).Value = 46081.52083333334; $wsDetailed.Cells.Item(75, 2).Value = 0.7; $wsDetailed.Cells.Item(75, 3).Value = "forecast"; $wsDetailed.Cells.Item(75, 4).Value = 46081; $wsDetailed.Cells.Item(75, 5).Value = "ON"
$wsDetailed.Cells.Item(76, 1).Value = 46081.54166666666; $wsDetailed.Cells.Item(76, 2).Value = 0.7; $wsDetailed.Cells.Item(76, 3).Value = "forecast"; $wsDetailed.Cells.Item(76, 4).Value = 46081; $wsDetailed.Cells.Item(76, 5).Value = "ON"
$wsDetailed.Cells.Item(77, 1).Value = 46081.5625; $wsDetailed.Cells.Item(77, 2).Value = -1.60912; $wsDetailed.Cells.Item(77, 3).Value = "forecast"; $wsDetailed.Cells.Item(77, 4).Value = 46081; $wsDetailed.Cells.Item(77, 5).Value = "ON"
$wsDetailed.Cells.Item(78, 1).Value = 46081.58333333334; $wsDetailed.Cells.Item(78, 2).Value = -4; $wsDetailed.Cells.Item(78, 3).Value = "forecast"; $wsDetailed.Cells.Item(78, 4).Value = 46081; $wsDetailed.Cells.Item(78, 5).Value = "ON"
$wsDetailed.Cells.Item(79, 1).Value = 46081.60416666666; $wsDetailed.Cells.Item(79, 2).Value = -3.75989; $wsDetailed.Cells.Item(79, 3).Value = "forecast"; $wsDetailed.Cells.Item(79, 4).Value = 46081; $wsDetailed.Cells.Item(79, 5).Value = "ON"
$wsDetailed.Cells.Item(80, 1).Value = 46081.625; $wsDetailed.Cells.Item(80, 2).Value = 0.51; $wsDetailed.Cells.Item(80, 3).Value = "forecast"; $wsDetailed.Cells.Item(80, 4).Value = 46081; $wsDetailed.Cells.Item(80, 5).Value = "ON"
$wsDetailed.Cells.Item(81, 1).Value = 46081.64583333334; $wsDetailed.Cells.Item(81, 2).Value = 35.88; $wsDetailed.Cells.Item(81, 3).Value = "forecast"; $wsDetailed.Cells.Item(81, 4).Value = 46081; $wsDetailed.Cells.Item(81, 5).Value = "ON"
$wsDetailed.Cells.Item(82, 1).Value = 46081.66666666666; $wsDetailed.Cells.Item(82, 2).Value = 9.64104; $wsDetailed.Cells.Item(82, 3).Value = "forecast"; $wsDetailed.Cells.Item(82, 4).Value = 46081; $wsDetailed.Cells.Item(82, 5).Value = "ON"
$wsDetailed.Cells.Item(83, 1).Value = 46081.6875; $wsDetailed.Cells.Item(83, 2).Value = -3.75989; $wsDetailed.Cells.Item(83, 3).Value = "forecast"; $wsDetailed.Cells.Item(83, 4).Value = 46081; $wsDetailed.Cells.Item(83, 5).Value = "ON"
$wsDetailed.Cells.Item(84, 1).Value = 46081.70833333334; $wsDetailed.Cells.Item(84, 2).Value = 35.88; $wsDetailed.Cells.Item(84, 3).Value = "forecast"; $wsDetailed.Cells.Item(84, 4).Value = 46081; $wsDetailed.Cells.Item(84, 5).Value = "ON"
$wsDetailed.Cells.Item(85, 1).Value = 46081.72916666666; $wsDetailed.Cells.Item(85, 2).Value = 50.88189; $wsDetailed.Cells.Item(85, 3).Value = "forecast"; $wsDetailed.Cells.Item(85, 4).Value = 46081; $wsDetailed.Cells.Item(85, 5).Value = "ON"
$wsDetailed.Cells.Item(86, 1).Value = 46081.75; $wsDetailed.Cells.Item(86, 2).Value = 50.55729; $wsDetailed.Cells.Item(86, 3).Value = "forecast"; $wsDetailed.Cells.Item(86, 4).Value = 46081; $wsDetailed.Cells.Item(86, 5).Value = "OFF"
$wsDetailed.Cells.Item(87, 1).Value = 46081.77083333334; $wsDetailed.Cells.Item(87, 2).Value = 65; $wsDetailed.Cells.Item(87, 3).Value = "forecast"; $wsDetailed.Cells.Item(87, 4).Value = 46081; $wsDetailed.Cells.Item(87, 5).Value = "OFF"
$wsDetailed.Cells.Item(88, 1).Value = 46081.79166666666; $wsDetailed.Cells.Item(88, 2).Value = 65.01002; $wsDetailed.Cells.Item(88, 3).Value = "forecast"; $wsDetailed.Cells.Item(88, 4).Value = 46081; $wsDetailed.Cells.Item(88, 5).Value = "OFF"
$wsDetailed.Cells.Item(89, 1).Value = 46081.8125; $wsDetailed.Cells.Item(89, 2).Value = 65.01003; $wsDetailed.Cells.Item(89, 3).Value = "forecast"; $wsDetailed.Cells.Item(89, 4).Value = 46081; $wsDetailed.Cells.Item(89, 5).Value = "OFF"
$wsDetailed.Cells.Item(90, 1).Value = 46081.83333333334; $wsDetailed.Cells.Item(90, 2).Value = 62.97923; $wsDetailed.Cells.Item(90, 3).Value = "forecast"; $wsDetailed.Cells.Item(90, 4).Value = 46081; $wsDetailed.Cells.Item(90, 5).Value = "OFF"
$wsDetailed.Cells.Item(91, 1).Value = 46081.85416666666; $wsDetailed.Cells.Item(91, 2).Value = 59.57788; $wsDetailed.Cells.Item(91, 3).Value = "forecast"; $wsDetailed.Cells.Item(91, 4).Value = 46081; $wsDetailed.Cells.Item(91, 5).Value = "OFF"
$wsDetailed.Cells.Item(92, 1).Value = 46081.875; $wsDetailed.Cells.Item(92, 2).Value = 57.32; $wsDetailed.Cells.Item(92, 3).Value = "forecast"; $wsDetailed.Cells.Item(92, 4).Value = 46081; $wsDetailed.Cells.Item(92, 5).Value = "OFF"
$wsDetailed.Cells.Item(93, 1).Value = 46081.89583333334; $wsDetailed.Cells.Item(93, 2).Value = 57.06; $wsDetailed.Cells.Item(93, 3).Value = "forecast"; $wsDetailed.Cells.Item(93, 4).Value = 46081; $wsDetailed.Cells.Item(93, 5).Value = "OFF"
$wsDetailed.Cells.Item(94, 1).Value = 46081.91666666666; $wsDetailed.Cells.Item(94, 2).Value = 56.98; $wsDetailed.Cells.Item(94, 3).Value = "forecast"; $wsDetailed.Cells.Item(94, 4).Value = 46081; $wsDetailed.Cells.Item(94, 5).Value = "OFF"
$wsDetailed.Cells.Item(95, 1).Value = 46081.9375; $wsDetailed.Cells.Item(95, 2).Value = 57.06; $wsDetailed.Cells.Item(95, 3).Value = "forecast"; $wsDetailed.Cells.Item(95, 4).Value = 46081; $wsDetailed.Cells.Item(95, 5).Value = "OFF"
$wsDetailed.Cells.Item(96, 1).Value = 46081.95833333334; $wsDetailed.Cells.Item(96, 2).Value = 57.06; $wsDetailed.Cells.Item(96, 3).Value = "forecast"; $wsDetailed.Cells.Item(96, 4).Value = 46081; $wsDetailed.Cells.Item(96, 5).Value = "OFF"
$wsDetailed.Cells.Item(97, 1).Value = 46081.97916666666; $wsDetailed.Cells.Item(97, 2).Value = 56.04; $wsDetailed.Cells.Item(97, 3).Value = "forecast"; $wsDetailed.Cells.Item(97, 4).Value = 46081; $wsDetailed.Cells.Item(97, 5).Value = "OFF"

# Number formats for the newly appended rows (A: datetime, D: date)
$wsDetailed.Range("A50:A97").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsDetailed.Range("D50:D97").NumberFormat = "YYYY-MM-DD"

Write-Host "edit applied"
